{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\n// New ranking text, in final top-to-bottom order, for the 14 surviving items.\nconst newTexts = [\n  \"\ud83e\udd47 Steak profit: 7.35 \u20aa/min\",\n  \"\ud83e\udd48 Stuffed Mushrooms profit: 6.81 \u20aa/min\",\n  \"\ud83e\udd49 Salad profit: 6.00 \u20aa/min\",\n  \"4) Brule Cream profit: 5.36 \u20aa/min\",\n  \"5) Pasta profit: 4.65 \u20aa/min\",\n  \"6) Pizza profit: 3.63 \u20aa/min\",\n  \"7) Krep profit: 3.56 \u20aa/min\",\n  \"8) Belgian Waffle profit: 3.20 \u20aa/min\",\n  \"9) Hamburger profit: 2.70 \u20aa/min\",\n  \"10) Empanadas profit: 2.68 \u20aa/min\",\n  \"11) Schnitzel profit: 2.48 \u20aa/min\",\n  \"12) Cake profit: 2.32 \u20aa/min\",\n  \"13) Roast profit: 2.00 \u20aa/min\",\n  \"14) Arancini profit: 1.81 \u20aa/min\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The first two paragraphs are the title and subtitle; the ranked list\n// starts at index 2. Match the existing \"profit:\" list paragraphs in\n// document order so the script is resilient to the exact starting index.\nconst listStart = paragraphs.items.findIndex((p) => p.text.includes(\"profit:\"));\n\nif (listStart === -1) {\n  throw new Error(\"Could not locate the profits list in the document body.\");\n}\n\nconst listParas = paragraphs.items.slice(listStart);\n\n// Update the surviving rows (renamed / renumbered / re-ranked) in place.\nfor (let i = 0; i < newTexts.length; i++) {\n  listParas[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\n\n// Remove the trailing rows that no longer exist in the updated ranking.\nfor (let i = newTexts.length; i < listParas.length; i++) {\n  listParas[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# New ranking text, in final top-to-bottom order, for the 14 surviving items.\n$newTexts = @(\n    \"\ud83e\udd47 Steak profit: 7.35 \u20aa/min\",\n    \"\ud83e\udd48 Stuffed Mushrooms profit: 6.81 \u20aa/min\",\n    \"\ud83e\udd49 Salad profit: 6.00 \u20aa/min\",\n    \"4) Brule Cream profit: 5.36 \u20aa/min\",\n    \"5) Pasta profit: 4.65 \u20aa/min\",\n    \"6) Pizza profit: 3.63 \u20aa/min\",\n    \"7) Krep profit: 3.56 \u20aa/min\",\n    \"8) Belgian Waffle profit: 3.20 \u20aa/min\",\n    \"9) Hamburger profit: 2.70 \u20aa/min\",\n    \"10) Empanadas profit: 2.68 \u20aa/min\",\n    \"11) Schnitzel profit: 2.48 \u20aa/min\",\n    \"12) Cake profit: 2.32 \u20aa/min\",\n    \"13) Roast profit: 2.00 \u20aa/min\",\n    \"14) Arancini profit: 1.81 \u20aa/min\"\n)\n\n# Find the first paragraph of the ranked \"profit:\" list.\n$listStart = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*profit:*\") {\n        $listStart = $i\n        break\n    }\n}\n\nif ($listStart -eq -1) {\n    throw \"Could not locate the profits list in the document body.\"\n}\n\n# Update the surviving rows (renamed / renumbered / re-ranked) in place.\nfor ($j = 0; $j -lt $newTexts.Count; $j++) {\n    $p = $d.Paragraphs.Item($listStart + $j)\n    $p.Range.Text = $newTexts[$j]\n}\n\n# Remove the trailing rows that no longer exist in the updated ranking,\n# starting from the end so indices of earlier paragraphs stay valid.\nfor ($i = $d.Paragraphs.Count; $i -ge $listStart + $newTexts.Count; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
